$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "69.312.48"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.746.28"
$ws.Range("E3").Value = "  +0.30%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "602.28"
$ws.Range("E5").Value = "  +0.08%  "
Set-TextValue $ws.Range("D6") "168.20"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "3.744.73"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +3.82%  "
$ws.Range("E11").Value = "  +0.13%  "
Set-TextValue $ws.Range("D12") "0.462"
$ws.Range("E12").Value = "  +0.71%  "
Set-TextValue $ws.Range("D13") "38.28"
$ws.Range("E13").Value = "  +0.70%  "
Set-TextValue $ws.Range("D14") "0.0000248"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "4.371.99"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "3.742.72"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "69.226.79"
$ws.Range("E17").Value = "  +0.72%  "
Set-TextValue $ws.Range("D18") "7.41"
$ws.Range("E18").Value = "  +2.37%  "
Set-TextValue $ws.Range("D19") "17.42"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  -1.46%  "
Set-TextValue $ws.Range("D21") "11.25"
$ws.Range("E21").Value = "  +12.20%  "
Set-TextValue $ws.Range("D22") "493.11"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +1.03%  "
Set-TextValue $ws.Range("D24") "0.0000150"
$ws.Range("E24").Value = "  +5.43%  "
Set-TextValue $ws.Range("D25") "84.85"
$ws.Range("E25").Value = "  -0.03%  "
Set-TextValue $ws.Range("D26") "2.30"
$ws.Range("E26").Value = "  -0.29%  "
Set-TextValue $ws.Range("D27") "12.33"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +0.01%  "
Set-TextValue $ws.Range("D30") "2.99"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("E32").Value = "  +1.04%  "
Set-TextValue $ws.Range("D33") "31.65"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "3.889.92"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "3.677.88"
$ws.Range("E36").Value = "  +0.25%  "
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("E39").Value = "  +5.42%  "
$ws.Range("E40").Value = "  -0.11%  "
Set-TextValue $ws.Range("D41") "0.328"
$ws.Range("E41").Value = "  +1.18%  "
Set-TextValue $ws.Range("D42") "3.06"
$ws.Range("E42").Value = "  +6.33%  "
Set-TextValue $ws.Range("D43") "48.92"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  +1.16%  "
Set-TextValue $ws.Range("D45") "425.62"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  -0.01%  "
Set-TextValue $ws.Range("D48") "40.24"
$ws.Range("E48").Value = "  -1.07%  "
Set-TextValue $ws.Range("D49") "141.41"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.788.34"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D51") "0.0356"
$ws.Range("E51").Value = "  +0.93%  "
